$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "NB" row (row 8); remaining rows shift up and SVM becomes row 8
$ws.Rows(8).Delete()

# Give the new header cells H1:L1 the same style as the existing header cells (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1:L1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row labels, B1:L1
$header = New-Object "object[,]" 1,11
$header[0,0] = "Algorithm"
$header[0,1] = "One Year Base mean"
$header[0,2] = "One Year Base std"
$header[0,3] = "Two Year Base mean"
$header[0,4] = "Two Year Base std"
$header[0,5] = "Three Year Base mean"
$header[0,6] = "Three Year Base std"
$header[0,7] = "Five Year Base mean"
$header[0,8] = "Five Year Base std"
$header[0,9] = "Ten Year Base mean"
$header[0,10] = "Ten Year Base std"
$ws.Range("B1:L1").Value = $header

# Data rows, A2:L8 (index, algorithm name, then 5 mean/std pairs)
$data = New-Object "object[,]" 7,12
$data[0,0] = 0
$data[0,1] = "LR"
$data[0,2] = 0.9119938560949127
$data[0,3] = 0.009790279344311809
$data[0,4] = 0.8963070843340659
$data[0,5] = 0.006677828304474398
$data[0,6] = 0.8822397829840913
$data[0,7] = 0.01181152650232513
$data[0,8] = 0.8724229535429913
$data[0,9] = 0.01944799787441679
$data[0,10] = 0.8582711231761113
$data[0,11] = 0.01339093123067229
$data[1,0] = 1
$data[1,1] = "LDA"
$data[1,2] = 0.9152466830857234
$data[1,3] = 0.008861120151360058
$data[1,4] = 0.9032202292741921
$data[1,5] = 0.008962786365398242
$data[1,6] = 0.8906384333704824
$data[1,7] = 0.01053375772981107
$data[1,8] = 0.8784767571747724
$data[1,9] = 0.02053679984776114
$data[1,10] = 0.8654100214907816
$data[1,11] = 0.01593204911929342
$data[2,0] = 2
$data[2,1] = "KNN"
$data[2,2] = 0.8984920950186701
$data[2,3] = 0.01022036510747025
$data[2,4] = 0.8919248698338074
$data[2,5] = 0.006618083869208864
$data[2,6] = 0.8881835327544609
$data[2,7] = 0.01347278070759444
$data[2,8] = 0.8909760411296327
$data[2,9] = 0.01436241370655155
$data[2,10] = 0.8830013573125213
$data[2,11] = 0.01738646189566103
$data[3,0] = 3
$data[3,1] = "DTREE"
$data[3,2] = 0.8861288101480363
$data[3,3] = 0.01370130374417349
$data[3,4] = 0.8848439425167924
$data[3,5] = 0.007922100351338578
$data[3,6] = 0.8769901289603566
$data[3,7] = 0.009519385145594259
$data[3,8] = 0.8820705877298505
$data[3,9] = 0.009479079206697282
$data[3,10] = 0.8863380839271576
$data[3,11] = 0.02013374691934713
$data[4,0] = 4
$data[4,1] = "RTREE"
$data[4,2] = 0.9111763459654141
$data[4,3] = 0.007060438046636296
$data[4,4] = 0.8974869549911709
$data[4,5] = 0.007767340313739636
$data[4,6] = 0.8880136063586151
$data[4,7] = 0.007750904569699226
$data[4,8] = 0.8750698144010999
$data[4,9] = 0.01729369143529494
$data[4,10] = 0.8637439203709988
$data[4,11] = 0.01597366520640282
$data[5,0] = 5
$data[5,1] = "XTREE"
$data[5,2] = 0.9147583485606843
$data[5,3] = 0.009489022536760284
$data[5,4] = 0.9070985288523232
$data[5,5] = 0.01200430918468063
$data[5,6] = 0.8993858155854653
$data[5,7] = 0.008849943255277649
$data[5,8] = 0.8936225439651716
$data[5,9] = 0.01392210896118182
$data[5,10] = 0.8939373374052708
$data[5,11] = 0.009508801975274454
$data[6,0] = 6
$data[6,1] = "SVM"
$data[6,2] = 0.9066248245544344
$data[6,3] = 0.008000900466325579
$data[6,4] = 0.9044015194099512
$data[6,5] = 0.008388464446453515
$data[6,6] = 0.8972854640980735
$data[6,7] = 0.01246652190714638
$data[6,8] = 0.8955168413816805
$data[6,9] = 0.01934739598232919
$data[6,10] = 0.8813392150209252
$data[6,11] = 0.01577020791687747
$ws.Range("A2:L8").Value = $data
